$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data per latest scrape

$ws.Range("D2").Value = "21.750.25"
$ws.Range("E2").Value = "  -1.43%  "

$ws.Range("E3").Value = "  -1.18%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("E5").Value = "  -0.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.88%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3880"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.60%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3195"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.08"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.22%  "

$ws.Range("E10").Value = "  -1.55%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.058"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.25%  "

$ws.Range("E12").Value = "  -0.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.641"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.601"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.67%  "

$ws.Range("D16").Value = "1.542.88"
$ws.Range("E16").Value = "  -1.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001109"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06587"
$ws.Range("D18").Style = "Normal"

$ws.Range("E19").Value = "  -2.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.145"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.84%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.42%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.384"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.69%  "

$ws.Range("D25").Value = "21.759.51"
$ws.Range("E25").Value = "  -1.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.382"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.56%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "146.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.81%  "

$ws.Range("E28").Value = "  -3.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.849"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.38%  "

$ws.Range("D30").Value = "1.717.25"
$ws.Range("E30").Value = "  -1.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.46%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.935"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.38%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9707"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -12.74%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08182"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.918"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.158"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06078"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.487"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -17.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02207"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2040"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.192"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.29%  "

$ws.Range("E42").Value = "  -0.10%  "

$ws.Range("E43").Value = "  -2.43%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5756"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.12%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.13%  "

$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.749"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.45%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5528"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.67%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.868"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.144"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06725"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.45%  "
